# Applies the "10 2-byte burst FEC + DW + Improved transmission method"
# timing-sample addition to Sheet1, per the commit:
#   "Improved WLESS_SendPacketBurst's transmission time to nearly 3 times faster."
#
# This mirrors the existing "Individual Packet Send Timing" burst blocks
# (rows 7/11/15 -> columns G/H/I) by adding a new sample in row 16, a new
# labeled block header in row 19 (G19/L19), and a new sample + derived
# "Frequency (Hz)" value in row 20 (G20/H20/I20/L20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New burst sample in row 16 (G/H start + end times, I = duration formula)
$ws.Range("G16").Value = 4.0006612500000003
$ws.Range("H16").Value = 4.0128336100000004
$ws.Range("I16").Formula = "=(H16-G16)*1000"

# New block label (row 19) for the improved transmission method, plus a new
# "Frequency (Hz)" column header next to it
$ws.Range("G19").Value = "10 2-byte burst FEC + DW + Improved transmission method"
$ws.Range("L19").Value = "Frequency (Hz)"

# New burst sample in row 20 for the improved method, with a derived
# frequency (1000 / duration-in-ms) in column L
$ws.Range("G20").Value = 4.00066364
$ws.Range("H20").Value = 4.0051203299999996
$ws.Range("I20").Formula = "=(H20-G20)*1000"
$ws.Range("L20").Formula = "=1000/I20"

# Widen the new column L to match the other data columns (target stored
# OOXML width is 22.140625; the host's ColumnWidth->stored-width quantizes
# to 1/6-character steps, so 21.33 is the closest achievable input)
$ws.Columns.Item(12).ColumnWidth = 21.33

# Move the active selection to the newly added frequency cell
$ws.Range("L21").Select()
